{"js": "// \"fixed HH get method\"\n// The \"Return a list of household IDs the user belongs to via a get request\"\n// bullet under \"Household Controller\" becomes\n// \"Return a list of households the user belongs to via a get request\".\n// The document's `_GoBack` bookmark (previously sitting at the very end of the\n// \"Registration Controller\" bullet above) also moves to sit right after the\n// word \"households\" in the edited sentence - this is simply where the\n// author's cursor ended up after making the last edit.\n\nconst body = context.document.body;\n\n// 1) Locate the paragraph that needs editing via a unique substring.\nconst target = body.search(\"household IDs the user belongs to\", { matchCase: false });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the Household IDs bullet, found \" + target.items.length);\n}\n\nconst paragraph = target.items[0].paragraphs.getFirst();\n\n// 2) Remove the old `_GoBack` bookmark from wherever it currently lives\n//    (the end of the \"...and return the user's id\" paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Turn \"household IDs\" into \"households\" by replacing \" IDs\" with \"s\".\nconst idsHits = paragraph.search(\" IDs\", { matchCase: false });\nidsHits.load(\"items\");\nawait context.sync();\n\nif (idsHits.items.length !== 1) {\n  throw new Error(\"Expected exactly one ' IDs' occurrence in the bullet, found \" + idsHits.items.length);\n}\n\nidsHits.items[0].insertText(\"s\", \"Replace\");\nawait context.sync();\n\n// 4) Re-insert the `_GoBack` bookmark immediately after \"households\".\nconst householdsHits = paragraph.search(\"households\", { matchCase: false });\nhouseholdsHits.load(\"items\");\nawait context.sync();\n\nif (householdsHits.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'households' occurrence in the bullet, found \" + householdsHits.items.length);\n}\n\nconst afterHouseholds = householdsHits.items[0].getRange(\"End\");\nafterHouseholds.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"fixed HH get method\"\n# The \"Return a list of household IDs the user belongs to via a get request\"\n# bullet under \"Household Controller\" becomes\n# \"Return a list of households the user belongs to via a get request\".\n# The document's `_GoBack` bookmark (previously sitting at the very end of the\n# \"Registration Controller\" bullet above) also moves to sit right after the\n# word \"households\" in the edited sentence - this is simply where the\n# author's cursor ended up after making the last edit.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old `_GoBack` bookmark from wherever it currently lives\n#    (the end of the \"...and return the user's id\" paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Turn \"household IDs\" into \"households\" in the target bullet.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"household IDs the user belongs to via a get request\")\nif (-not $found) {\n    throw \"Could not find the 'household IDs ...' bullet to update\"\n}\n$rng.Text = \"households the user belongs to via a get request\"\n\n# 3) Re-insert the `_GoBack` bookmark immediately after \"households\" in the\n#    edited sentence. Search on a phrase that is unique to this bullet (other\n#    bullets elsewhere in the doc also start with \"Return a list of\n#    households...\") and then narrow the found range down to just the word.\n$hhRange = $d.Content\n$hhFound = $hhRange.Find.Execute(\"households the user belongs to\")\nif (-not $hhFound) {\n    throw \"Could not find 'households the user belongs to' to anchor the bookmark\"\n}\n$hhRange.Collapse(1)       # wdCollapseStart - start of the matched phrase\n$hhRange.MoveEnd(1, 10)    # wdCharacter, Len(\"households\") -> select just the word\n$hhRange.Collapse(0)       # wdCollapseEnd - collapse to right after \"households\"\n$d.Bookmarks.Add(\"_GoBack\", $hhRange)\n"}
